$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on the cells we are about to update so that the
# numeric-looking / percentage-looking strings are stored as literal text
# (matching the original inline-string cell contents) instead of being
# auto-converted to numbers by Excel.
$textCells = @("D2","E2","D3","E3","D4","E4","D5","E5","D6","E6","D7","E7","D8","E8","D9","E9","D10","E10","D11","E11","D12","E12","D13","E13","D14","E14","D15","E15","E16","D17","E17","D18","E18","D19","E19","D20","E20","D21","E21","D22","E22","D23","E23","D24","E24","D25","E25","D26","E26","D27","D39","E39","D40","E40","D41","E41","D42","E42","D43","E43","D44","E44","D45","E45","D46","E46","D47","E47","D48","E48","D49","E49","D50","E50","D51","E51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated values
$ws.Range("D2").Value = "308.74"
$ws.Range("E2").Value = "-1.02%"
$ws.Range("D3").Value = "36.30"
$ws.Range("E3").Value = "-3.79%"
$ws.Range("D4").Value = "5.117"
$ws.Range("E4").Value = "-0.28%"
$ws.Range("D5").Value = "0.07716"
$ws.Range("E5").Value = "-2.33%"
$ws.Range("D6").Value = "4.385"
$ws.Range("E6").Value = "-0.61%"
$ws.Range("D7").Value = "8.297"
$ws.Range("E7").Value = "0.24%"
$ws.Range("D8").Value = "1.858"
$ws.Range("E8").Value = "-2.43%"
$ws.Range("D9").Value = "2.942"
$ws.Range("E9").Value = "-5.07%"
$ws.Range("D10").Value = "0.9192"
$ws.Range("E10").Value = "-0.70%"
$ws.Range("D11").Value = "0.1131"
$ws.Range("E11").Value = "-5.95%"
$ws.Range("D12").Value = "0.1854"
$ws.Range("E12").Value = "-3.97%"
$ws.Range("D13").Value = "0.08769"
$ws.Range("E13").Value = "-3.75%"
$ws.Range("D14").Value = "0.03324"
$ws.Range("E14").Value = "-0.05%"
$ws.Range("D15").Value = "0.09529"
$ws.Range("E15").Value = "-1.04%"
$ws.Range("E16").Value = "-0.39%"
$ws.Range("D17").Value = "0.006144"
$ws.Range("E17").Value = "4.59%"
$ws.Range("D18").Value = "3.363"
$ws.Range("E18").Value = "-4.28%"
$ws.Range("D19").Value = "0.3446"
$ws.Range("E19").Value = "1.36%"
$ws.Range("D20").Value = "6.327"
$ws.Range("E20").Value = "19.65%"
$ws.Range("D21").Value = "0.1316"
$ws.Range("E21").Value = "3.44%"
$ws.Range("D22").Value = "0.2314"
$ws.Range("E22").Value = "-10.63%"
$ws.Range("D23").Value = "0.04337"
$ws.Range("E23").Value = "-0.84%"
$ws.Range("D24").Value = "0.001203"
$ws.Range("E24").Value = "-3.63%"
$ws.Range("D25").Value = "0.004261"
$ws.Range("E25").Value = "-1.17%"
$ws.Range("D26").Value = "0.0001331"
$ws.Range("E26").Value = "9.05%"
$ws.Range("D27").Value = "0.0002904"
$ws.Range("D39").Value = "0.02102"
$ws.Range("E39").Value = "-0.70%"
$ws.Range("D40").Value = "0.04921"
$ws.Range("E40").Value = "-5.07%"
$ws.Range("D41").Value = "0.007567"
$ws.Range("E41").Value = "-1.31%"
$ws.Range("D42").Value = "0.1349"
$ws.Range("E42").Value = "-0.89%"
$ws.Range("D43").Value = "0.008554"
$ws.Range("E43").Value = "-5.47%"
$ws.Range("D44").Value = "0.002071"
$ws.Range("E44").Value = "3.01%"
$ws.Range("D45").Value = "0.008383"
$ws.Range("E45").Value = "-2.58%"
$ws.Range("D46").Value = "0.00006469"
$ws.Range("E46").Value = "-3.44%"
$ws.Range("D47").Value = "0.00000000751"
$ws.Range("E47").Value = "0.14%"
$ws.Range("D48").Value = "0.003297"
$ws.Range("E48").Value = "18.25%"
$ws.Range("D49").Value = "0.001445"
$ws.Range("E49").Value = "20.47%"
$ws.Range("D50").Value = "0.00002103"
$ws.Range("E50").Value = "0.14%"
$ws.Range("D51").Value = "0.0002003"
$ws.Range("E51").Value = "0.14%"
